# Realestate Update resale numbers 2024-01-06 22:00
# Append a new data row (row 27) to the CityResaleNum sheet with the
# latest resale-number snapshot for 2024-01-06 22:00:34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

# Column A (date) and D (zero-padded week number) look numeric/date-like
# to Excel's auto-detection, so force literal text storage for those two
# ("2024-01-06" would otherwise become a date serial, "00" would become 0).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-06"
$ws.Cells.Item($row, 2).Value = "22:00:34"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "00"

# Columns E-T are the per-city resale-number counts (plain numbers; -1
# marks a city with no data for this snapshot).
$ws.Cells.Item($row, 5).Value = 140580
$ws.Cells.Item($row, 6).Value = 142969
$ws.Cells.Item($row, 7).Value = 172338
$ws.Cells.Item($row, 8).Value = 147281
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118372
$ws.Cells.Item($row, 11).Value = 224616
$ws.Cells.Item($row, 12).Value = 249355
$ws.Cells.Item($row, 13).Value = 185131
$ws.Cells.Item($row, 14).Value = 110396
$ws.Cells.Item($row, 15).Value = 40635
$ws.Cells.Item($row, 16).Value = 30809
$ws.Cells.Item($row, 17).Value = 72509
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42264
$ws.Cells.Item($row, 20).Value = -1
